$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "duplicate" status values in column B for rows 2, 3 and 5
$ws.Range("B2").Value = $null
$ws.Range("B3").Value = $null
$ws.Range("B5").Value = $null

# Update the selection to match the new state
$ws.Range("B3:B5").Select()
